$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted above the existing row 373,
# pushing the old rows 373-454 down to 374-455.
$ws.Rows.Item(373).Insert()

# Fill in the new record's values.
$ws.Range("A373").Value = 5
$ws.Range("B373").Value = "Macroferia Regional de Talca"
$ws.Range("C373").Value = "Maule"
$ws.Range("D373").Value = 44889
$ws.Range("E373").Value = 7
$ws.Range("F373").Value = 100112032
$ws.Range("G373").Value = "Zapallo italiano"
$ws.Range("H373").Value = "Sin especificar"
$ws.Range("I373").Value = "Primera"
$ws.Range("J373").Value = 400
$ws.Range("K373").Value = 6000
$ws.Range("L373").Value = 6000
$ws.Range("M373").Value = 6000
$ws.Range("N373").Value = "$/caja 50 unidades"
$ws.Range("O373").Value = "Región del Maule"
$ws.Range("P373").Value = 120
$ws.Range("Q373").Value = 50
$ws.Range("R373").Value = "Hortaliza"
